$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 - Brooklyn Nets
$ws.Range("Y4").Value = 3.9
$ws.Range("AN4").Value = 2

# Row 7 - Cleveland Cavaliers
$ws.Range("AX7").Value = 17

# Row 8 - Dallas Mavericks
$ws.Range("AY8").Value = 6

# Row 10 - Detroit Pistons
$ws.Range("S10").Value = 31.7
$ws.Range("T10").Value = 42.8
$ws.Range("V10").Value = 14
$ws.Range("AV10").Value = 9

# Row 14 - LA Clippers
$ws.Range("AX14").Value = 19

# Row 16 - Memphis Grizzlies
$ws.Range("J16").Value = 92.2
$ws.Range("R16").Value = 10.6
$ws.Range("T16").Value = 45.4
$ws.Range("AR16").Value = 12
$ws.Range("AX16").Value = 17

# Row 18 - Milwaukee Bucks
$ws.Range("M18").Value = 39.6
$ws.Range("X18").Value = 4.8
$ws.Range("AN18").Value = 3
$ws.Range("AX18").Value = 19

# Row 21 - New York Knicks
$ws.Range("AR21").Value = 13

# Row 27 - Sacramento Kings
$ws.Range("AV27").Value = 8
